$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Add($ws.Cells.Item(2,1), "https://example.com")
